# Update the three due_date values (C2:C4) from *-01 to *-15.
# These cells hold plain text dates (e.g. "2024-01-01") stored as shared
# strings, not real Excel dates. A straight `.Value = "2024-01-15"` would be
# auto-parsed by Excel into a date serial number, which we don't want - the
# target keeps them as literal text. Prefixing with a leading apostrophe
# forces text entry (exactly like typing '2024-01-15 into a cell), and the
# trailing ClearFormats() strips the "treat as text" formatting residue so
# the cells end up with no explicit style, just like the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'2024-01-15"
$ws.Range("C3").Value = "'2024-02-15"
$ws.Range("C4").Value = "'2024-03-15"

$ws.Range("C2:C4").ClearFormats()
